$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.239.80"
$ws.Range("E2").Value = "  +2.30%  "

$ws.Range("D3").Value = "2.347.27"
$ws.Range("E3").Value = "  +6.19%  "

$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.35"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.18%  "

$ws.Range("E7").Value = "  +3.59%  "

$ws.Range("E8").Value = "  -0.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.639"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +7.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.00"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0938"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.83"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.46%  "

$ws.Range("E13").Value = "  +9.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.106"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.34"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +9.78%  "

$ws.Range("D16").Value = "2.701.57"
$ws.Range("E16").Value = "  +6.04%  "

$ws.Range("D17").Value = "2.349.59"
$ws.Range("E17").Value = "  +5.09%  "

$ws.Range("D18").Value = "43.207.27"
$ws.Range("E18").Value = "  +2.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000109"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.25"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.39"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.03%  "

$ws.Range("E22").Value = "  +12.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.44"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "255.18"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +12.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.09"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.02"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.32%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.11"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.37%  "

$ws.Range("E29").Value = "  +1.06%  "

$ws.Range("E30").Value = "  +7.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.69"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("E32").Value = "  -0.73%  "

$ws.Range("E33").Value = "  +4.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.08"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +10.81%  "

$ws.Range("E35").Value = "  +5.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.96"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.97%  "

$ws.Range("E37").Value = "  +3.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.11"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.105"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.64%  "

$ws.Range("E40").Value = "  +11.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.65"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.93%  "

$ws.Range("E42").Value = "  +14.70%  "

$ws.Range("E43").Value = "  +1.73%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.70"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.69%  "

$ws.Range("E46").Value = "  +4.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.28"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +12.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.72"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +7.44%  "

$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("E50").Value = "  +3.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "69.95"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.74%  "
